# The first match row (row 2, " Oct 29 2020" vs Chennai Super Kings) is
# removed; the second match row (the old row 3, " Oct 7 2020") shifts up
# to become the new row 2, and the sheet's used range shrinks by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2").Delete()
